# "Se agrega funcionalidad de CITI ventas"
#
# - Mark the two "citi ventas" backlog rows (113 & 114) as finished.
# - Add a new backlog row (116) for the new quotation-report task.
# - Update the sheet selection / scroll position to match where the
#   user ended up working (row ~106..116 area).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")
$ws.Activate()

# citi ventas rows -> "terminado"
$ws.Range("B113").Value = "terminado"
$ws.Range("B114").Value = "terminado"

# New backlog item: quotation report on the articles table.
$ws.Range("A116").Value = "agregar reporte de cotizaion a tabla de articulos"
$ws.Range("B116").Value = "no comenzado"

# Reflect the new used range / scroll & selection state.
$win = $excel.ActiveWindow
$win.ScrollRow = 48
$win.ScrollColumn = 1
$win.Height = 7755

[void]$ws.Range("C106").Select()
